$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1; Col=1; Text="160÷4=40, 0"},
    @{Row=1; Col=2; Text="629÷8=78, 5"},
    @{Row=1; Col=3; Text="991÷6=165, 1"},
    @{Row=1; Col=4; Text="890÷8=111, 2"},
    @{Row=1; Col=5; Text="491÷4=122, 3"},

    @{Row=5; Col=1; Text="258÷4=64, 2"},
    @{Row=5; Col=2; Text="334÷4=83, 2"},
    @{Row=5; Col=3; Text="908÷5=181, 3"},
    @{Row=5; Col=4; Text="652÷2=326, 0"},
    @{Row=5; Col=5; Text="337÷2=168, 1"},

    @{Row=9; Col=1; Text="528÷9=58, 6"},
    @{Row=9; Col=2; Text="483÷8=60, 3"},
    @{Row=9; Col=3; Text="545÷3=181, 2"},
    @{Row=9; Col=4; Text="969÷6=161, 3"},
    @{Row=9; Col=5; Text="231÷8=28, 7"},

    @{Row=13; Col=1; Text="378÷2=189, 0"},
    @{Row=13; Col=2; Text="986÷2=493, 0"},
    @{Row=13; Col=3; Text="880÷8=110, 0"},
    @{Row=13; Col=4; Text="256÷3=85, 1"},
    @{Row=13; Col=5; Text="404÷7=57, 5"},

    @{Row=17; Col=1; Text="578÷2=289, 0"},
    @{Row=17; Col=2; Text="330÷8=41, 2"},
    @{Row=17; Col=3; Text="779÷8=97, 3"},
    @{Row=17; Col=4; Text="731÷8=91, 3"},
    @{Row=17; Col=5; Text="298÷8=37, 2"}
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    $cell.Range.Text = $item.Text
}

Write-Host "Done applying replacements."
